# The underlying source record for row 2 (artfynd id 465378 / Loweomyces
# wynneae, "Lundticka") was replaced by a different observation record
# (id 105278762 / Epipactis helleborine, "Skogsknipprot") that, in the
# original worksheet, lived in row 5 - and vice versa: row 5 now carries
# what used to be row 2's record. Net effect: the full data rows for rows 2
# and 5 are swapped, column by column (A:AY), including clearing cells that
# are populated on one side but not the other.
#
# Cells are written through `.Formula = "'<text>"` (a literal leading
# single quote = Excel's "treat as text" prefix) rather than `.Value = <text>`
# so that numeric-looking strings ("20", "4") and date-looking strings
# ("2022-08-18") are preserved as text instead of being auto-converted to
# numbers/dates by Excel - matching how this sheet stores them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: adopt what used to be row 5's values ---
    $ws.Range("A2").NumberFormat = "General"
    $ws.Range("A2").Value = 105278762
    $ws.Range("B2").NumberFormat = "General"
    $ws.Range("B2").Value = 96312
    $ws.Range("C2").Formula = "'Ovaliderad"
    $ws.Range("D2").Formula = "'LC"
    $ws.Range("E2").NumberFormat = "General"
    $ws.Range("E2").Value = 219798
    $ws.Range("F2").Formula = "'Skogsknipprot"
    $ws.Range("G2").Formula = "'Epipactis helleborine"
    $ws.Range("H2").Formula = "'(L.) Crantz"
    $ws.Range("I2").Formula = "'4"
    $ws.Range("J2").ClearContents()
    $ws.Range("K2").ClearContents()
    $ws.Range("P2").Formula = "'Trolleholm, Sk"
    $ws.Range("Q2").NumberFormat = "General"
    $ws.Range("Q2").Value = 392086.6705776053
    $ws.Range("R2").NumberFormat = "General"
    $ws.Range("R2").Value = 6199047.849894262
    $ws.Range("S2").NumberFormat = "General"
    $ws.Range("S2").Value = 10
    $ws.Range("T2").Formula = "'Skåne"
    $ws.Range("U2").Formula = "'Svalöv"
    $ws.Range("V2").Formula = "'Skåne"
    $ws.Range("W2").Formula = "'Torrlösa"
    $ws.Range("Y2").Formula = "'2022-08-18"
    $ws.Range("Z2").Formula = "'00:00"
    $ws.Range("AA2").Formula = "'2022-09-26"
    $ws.Range("AB2").Formula = "'00:00"
    $ws.Range("AC2").ClearContents()
    $ws.Range("AD2").NumberFormat = "General"
    $ws.Range("AD2").Value = $false
    $ws.Range("AE2").NumberFormat = "General"
    $ws.Range("AE2").Value = $false
    $ws.Range("AG2").NumberFormat = "General"
    $ws.Range("AG2").Value = $false
    $ws.Range("AH2").ClearContents()
    $ws.Range("AI2").ClearContents()
    $ws.Range("AQ2").ClearContents()
    $ws.Range("AR2").ClearContents()
    $ws.Range("AT2").Formula = "'"
    $ws.Range("AW2").Formula = "'Örjan Fritz"
    $ws.Range("AX2").Formula = "'Örjan Fritz"
    $ws.Range("AY2").Formula = "'"

# --- Row 5: adopt what used to be row 2's values ---
    $ws.Range("A5").NumberFormat = "General"
    $ws.Range("A5").Value = 465378
    $ws.Range("B5").NumberFormat = "General"
    $ws.Range("B5").Value = 90042
    $ws.Range("C5").Formula = "'Ovaliderad"
    $ws.Range("D5").Formula = "'VU"
    $ws.Range("E5").NumberFormat = "General"
    $ws.Range("E5").Value = 1627
    $ws.Range("F5").Formula = "'Lundticka"
    $ws.Range("G5").Formula = "'Loweomyces wynneae"
    $ws.Range("H5").Formula = "'(Berk. & Broome) Jülich"
    $ws.Range("I5").Formula = "'20"
    $ws.Range("J5").Formula = "'fruktkroppar"
    $ws.Range("K5").Formula = "'frukt-/fröspridning"
    $ws.Range("P5").Formula = "'Trolleholms gods, Sk"
    $ws.Range("Q5").NumberFormat = "General"
    $ws.Range("Q5").Value = 392454.3760533207
    $ws.Range("R5").NumberFormat = "General"
    $ws.Range("R5").Value = 6198116.700040066
    $ws.Range("S5").NumberFormat = "General"
    $ws.Range("S5").Value = 50
    $ws.Range("T5").Formula = "'Skåne"
    $ws.Range("U5").Formula = "'Svalöv"
    $ws.Range("V5").Formula = "'Skåne"
    $ws.Range("W5").Formula = "'Torrlösa"
    $ws.Range("Y5").Formula = "'2008-11-17"
    $ws.Range("Z5").Formula = "'00:00"
    $ws.Range("AA5").Formula = "'2008-11-17"
    $ws.Range("AB5").Formula = "'00:00"
    $ws.Range("AC5").Formula = "'Fruktkropparna i 2 fläckar åtskilja med ca 3 meters mellanrum"
    $ws.Range("AD5").NumberFormat = "General"
    $ws.Range("AD5").Value = $false
    $ws.Range("AE5").NumberFormat = "General"
    $ws.Range("AE5").Value = $false
    $ws.Range("AG5").NumberFormat = "General"
    $ws.Range("AG5").Value = $false
    $ws.Range("AH5").Formula = "'Ädellövskog"
    $ws.Range("AI5").Formula = "'Ask-al-skog"
    $ws.Range("AQ5").Formula = "'Sten Svantesson"
    $ws.Range("AR5").Formula = "'"
    $ws.Range("AT5").Formula = "'"
    $ws.Range("AW5").Formula = "'Sten Svantesson"
    $ws.Range("AX5").Formula = "'Sten Svantesson"
    $ws.Range("AY5").Formula = "'"
